$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 1596.2354
$ws.Cells.Item(19, 9).Value = 1576.125
$ws.Cells.Item(19, 10).Value = 1614.1111
$ws.Cells.Item(19, 11).Value = 1576.125
$ws.Cells.Item(19, 12).Value = 1614.1111
$ws.Cells.Item(19, 13).Value = -1401.125
$ws.Cells.Item(19, 14).Value = -1964.1111

$ws.Cells.Item(43, 8).Value = 5800
$ws.Cells.Item(43, 9).Value = 5800
$ws.Cells.Item(43, 11).Value = 5800
$ws.Cells.Item(43, 13).Value = -5731

$ws.Cells.Item(74, 8).Value = 7193.5557
$ws.Cells.Item(74, 9).Value = 6842.875
$ws.Cells.Item(74, 11).Value = 6842.875
$ws.Cells.Item(74, 13).Value = -5906.875

$ws.Cells.Item(76, 8).Value = 5050
$ws.Cells.Item(76, 9).Value = 4000
$ws.Cells.Item(76, 10).Value = 5260
$ws.Cells.Item(76, 11).Value = 4000
$ws.Cells.Item(76, 12).Value = 5260
$ws.Cells.Item(76, 13).Value = -3685
$ws.Cells.Item(76, 14).Value = -5890

$ws.Cells.Item(77, 8).Value = 7193.5557
$ws.Cells.Item(77, 9).Value = 6842.875
$ws.Cells.Item(77, 11).Value = 34214.375
$ws.Cells.Item(77, 13).Value = -29534.375

$ws.Cells.Item(79, 8).Value = 5050
$ws.Cells.Item(79, 9).Value = 4000
$ws.Cells.Item(79, 10).Value = 5260
$ws.Cells.Item(79, 11).Value = 4000
$ws.Cells.Item(79, 12).Value = 5260
$ws.Cells.Item(79, 13).Value = -2908
$ws.Cells.Item(79, 14).Value = -7444

$ws.Cells.Item(86, 8).Value = 1612.3334
$ws.Cells.Item(86, 9).Value = 1425
$ws.Cells.Item(86, 10).Value = 2549
$ws.Cells.Item(86, 11).Value = 1425
$ws.Cells.Item(86, 12).Value = 2549
$ws.Cells.Item(86, 13).Value = -302
$ws.Cells.Item(86, 14).Value = -4795

$ws.Cells.Item(89, 8).Value = 1612.3334
$ws.Cells.Item(89, 9).Value = 1425
$ws.Cells.Item(89, 10).Value = 2549
$ws.Cells.Item(89, 11).Value = 7125
$ws.Cells.Item(89, 12).Value = 12745
$ws.Cells.Item(89, 13).Value = -1509
$ws.Cells.Item(89, 14).Value = -23977

$ws.Cells.Item(137, 8).Value = 21160.25
$ws.Cells.Item(137, 9).Value = 17522.7
$ws.Cells.Item(137, 10).Value = 27222.834
$ws.Cells.Item(137, 11).Value = 52568.10000000001
$ws.Cells.Item(137, 12).Value = 81668.50199999999
$ws.Cells.Item(137, 13).Value = -50018.10000000001
$ws.Cells.Item(137, 14).Value = -86768.50199999999

$ws.Cells.Item(138, 8).Value = 44345.383
$ws.Cells.Item(138, 9).Value = 4266.8
$ws.Cells.Item(138, 10).Value = 98998
$ws.Cells.Item(138, 11).Value = 12800.4
$ws.Cells.Item(138, 12).Value = 296994
$ws.Cells.Item(138, 13).Value = -7660.400000000001
$ws.Cells.Item(138, 14).Value = -307274

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 14816.4
$ws.Cells.Item(32, 9).Value = 15024.513
$ws.Cells.Item(32, 11).Value = 15024.513
$ws.Cells.Item(32, 13).Value = -14737.513

$ws.Cells.Item(61, 8).Value = 9908.166999999999
$ws.Cells.Item(61, 9).Value = 4444.15
$ws.Cells.Item(61, 10).Value = 37228.25
$ws.Cells.Item(61, 11).Value = 4444.15
$ws.Cells.Item(61, 12).Value = 37228.25
$ws.Cells.Item(61, 13).Value = -4232.15
$ws.Cells.Item(61, 14).Value = -37652.25

$ws.Cells.Item(136, 8).Value = 9908.166999999999
$ws.Cells.Item(136, 9).Value = 4444.15
$ws.Cells.Item(136, 10).Value = 37228.25
$ws.Cells.Item(136, 11).Value = 13332.45
$ws.Cells.Item(136, 12).Value = 111684.75
$ws.Cells.Item(136, 13).Value = -10782.45
$ws.Cells.Item(136, 14).Value = -116784.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 1330.963
$ws.Cells.Item(99, 9).Value = 1001
$ws.Cells.Item(99, 10).Value = 2114.625
$ws.Cells.Item(99, 11).Value = 1001
$ws.Cells.Item(99, 12).Value = 2114.625
$ws.Cells.Item(99, 13).Value = 497
$ws.Cells.Item(99, 14).Value = -5110.625

$ws.Cells.Item(107, 8).Value = 2519.5
$ws.Cells.Item(107, 9).Value = 2519.5
$ws.Cells.Item(107, 11).Value = 2519.5
$ws.Cells.Item(107, 13).Value = -599.5

$ws.Cells.Item(134, 8).Value = 3965.5151
$ws.Cells.Item(134, 9).Value = 3933.1875
$ws.Cells.Item(134, 10).Value = 5000
$ws.Cells.Item(134, 11).Value = 11799.5625
$ws.Cells.Item(134, 12).Value = 15000
$ws.Cells.Item(134, 13).Value = -9264.5625
$ws.Cells.Item(134, 14).Value = -20070

$ws.Cells.Item(137, 8).Value = 67399.60000000001
$ws.Cells.Item(137, 10).Value = 63750
$ws.Cells.Item(137, 12).Value = 63750
$ws.Cells.Item(137, 14).Value = -73950

$ws.Cells.Item(138, 8).Value = 65000
$ws.Cells.Item(138, 10).Value = 65000
$ws.Cells.Item(138, 12).Value = 65000
$ws.Cells.Item(138, 14).Value = -75280

$ws.Cells.Item(139, 8).Value = 109999.75
$ws.Cells.Item(139, 10).Value = 109999.75
$ws.Cells.Item(139, 12).Value = 109999.75
$ws.Cells.Item(139, 14).Value = -120279.75

$ws.Cells.Item(140, 8).Value = 103798.6
$ws.Cells.Item(140, 10).Value = 103798.6
$ws.Cells.Item(140, 12).Value = 103798.6
$ws.Cells.Item(140, 14).Value = -114158.6

$ws.Cells.Item(141, 8).Value = 38500
$ws.Cells.Item(141, 9).Value = 38500
$ws.Cells.Item(141, 11).Value = 38500
$ws.Cells.Item(141, 13).Value = -33320

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1451162.8
$ws.Cells.Item(31, 9).Value = 1962539.5
$ws.Cells.Item(31, 10).Value = 2261.8333
$ws.Cells.Item(31, 11).Value = 1962539.5
$ws.Cells.Item(31, 12).Value = 2261.8333
$ws.Cells.Item(31, 13).Value = -1962244.5
$ws.Cells.Item(31, 14).Value = -2851.8333

$ws.Cells.Item(34, 8).Value = 1451162.8
$ws.Cells.Item(34, 9).Value = 1962539.5
$ws.Cells.Item(34, 10).Value = 2261.8333
$ws.Cells.Item(34, 11).Value = 1962539.5
$ws.Cells.Item(34, 12).Value = 2261.8333
$ws.Cells.Item(34, 13).Value = -1962337.5
$ws.Cells.Item(34, 14).Value = -2665.8333

$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 11).Value = 0
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 13).ClearContents()
$ws.Cells.Item(74, 14).ClearContents()

$ws.Cells.Item(77, 8).Value = 0
$ws.Cells.Item(77, 9).Value = 0
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 11).Value = 0
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 13).ClearContents()
$ws.Cells.Item(77, 14).ClearContents()

$ws.Cells.Item(86, 8).Value = 12491.737
$ws.Cells.Item(86, 9).Value = 11805.091
$ws.Cells.Item(86, 10).Value = 13435.875
$ws.Cells.Item(86, 11).Value = 11805.091
$ws.Cells.Item(86, 12).Value = 13435.875
$ws.Cells.Item(86, 13).Value = -10682.091
$ws.Cells.Item(86, 14).Value = -15681.875

$ws.Cells.Item(89, 8).Value = 12491.737
$ws.Cells.Item(89, 9).Value = 11805.091
$ws.Cells.Item(89, 10).Value = 13435.875
$ws.Cells.Item(89, 11).Value = 59025.455
$ws.Cells.Item(89, 12).Value = 67179.375
$ws.Cells.Item(89, 13).Value = -53409.455
$ws.Cells.Item(89, 14).Value = -78411.375

$ws.Cells.Item(99, 8).Value = 7158.5
$ws.Cells.Item(99, 9).Value = 6665
$ws.Cells.Item(99, 10).Value = 7898.75
$ws.Cells.Item(99, 11).Value = 6665
$ws.Cells.Item(99, 12).Value = 7898.75
$ws.Cells.Item(99, 13).Value = -5167
$ws.Cells.Item(99, 14).Value = -10894.75

$ws.Cells.Item(107, 8).Value = 625.9524
$ws.Cells.Item(107, 9).Value = 443.88235
$ws.Cells.Item(107, 11).Value = 443.88235
$ws.Cells.Item(107, 13).Value = 1476.11765

$ws.Cells.Item(122, 8).Value = 921.125
$ws.Cells.Item(122, 9).Value = 614.6
$ws.Cells.Item(122, 10).Value = 1432
$ws.Cells.Item(122, 11).Value = 1843.8
$ws.Cells.Item(122, 12).Value = 4296
$ws.Cells.Item(122, 13).Value = 606.1999999999998
$ws.Cells.Item(122, 14).Value = -9196

$ws.Cells.Item(126, 8).Value = 7158.5
$ws.Cells.Item(126, 9).Value = 6665
$ws.Cells.Item(126, 10).Value = 7898.75
$ws.Cells.Item(126, 11).Value = 19995
$ws.Cells.Item(126, 12).Value = 23696.25
$ws.Cells.Item(126, 13).Value = -17525
$ws.Cells.Item(126, 14).Value = -28636.25

$ws.Cells.Item(132, 8).Value = 44510.695
$ws.Cells.Item(132, 9).Value = 56172.445
$ws.Cells.Item(132, 10).Value = 2528.4
$ws.Cells.Item(132, 11).Value = 168517.335
$ws.Cells.Item(132, 12).Value = 7585.200000000001
$ws.Cells.Item(132, 13).Value = -165987.335
$ws.Cells.Item(132, 14).Value = -12645.2

$ws.Cells.Item(134, 8).Value = 2027.7778
$ws.Cells.Item(134, 9).Value = 1710
$ws.Cells.Item(134, 10).Value = 6000
$ws.Cells.Item(134, 11).Value = 5130
$ws.Cells.Item(134, 12).Value = 18000
$ws.Cells.Item(134, 13).Value = -2595
$ws.Cells.Item(134, 14).Value = -23070

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 876.75
$ws.Cells.Item(113, 9).Value = 534.875
$ws.Cells.Item(113, 10).Value = 1047.6875
$ws.Cells.Item(113, 11).Value = 1604.625
$ws.Cells.Item(113, 12).Value = 3143.0625
$ws.Cells.Item(113, 13).Value = 565.375
$ws.Cells.Item(113, 14).Value = -7483.0625

$ws.Cells.Item(132, 8).Value = 1336.8
$ws.Cells.Item(132, 9).Value = 800
$ws.Cells.Item(132, 10).Value = 1471
$ws.Cells.Item(132, 11).Value = 7200
$ws.Cells.Item(132, 12).Value = 13239
$ws.Cells.Item(132, 13).Value = -4670
$ws.Cells.Item(132, 14).Value = -18299

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(98, 8).Value = 14752.5
$ws.Cells.Item(98, 10).Value = 14752.5
$ws.Cells.Item(98, 12).Value = 14752.5
$ws.Cells.Item(98, 14).Value = -20742.5

$ws.Cells.Item(102, 8).Value = 35127.465
$ws.Cells.Item(102, 9).Value = 37279.43
$ws.Cells.Item(102, 10).Value = 5000
$ws.Cells.Item(102, 11).Value = 37279.43
$ws.Cells.Item(102, 12).Value = 5000
$ws.Cells.Item(102, 13).Value = -35657.43
$ws.Cells.Item(102, 14).Value = -8244

$ws.Cells.Item(113, 8).Value = 1308.6
$ws.Cells.Item(113, 9).Value = 1223.875
$ws.Cells.Item(113, 11).Value = 1223.875
$ws.Cells.Item(113, 13).Value = 946.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 1599.4736
$ws.Cells.Item(61, 9).Value = 1292.7333
$ws.Cells.Item(61, 10).Value = 2749.75
$ws.Cells.Item(61, 11).Value = 1292.7333
$ws.Cells.Item(61, 12).Value = 2749.75
$ws.Cells.Item(61, 13).Value = -1090.7333
$ws.Cells.Item(61, 14).Value = -3153.75

$ws.Cells.Item(113, 8).Value = 1599.4736
$ws.Cells.Item(113, 9).Value = 1292.7333
$ws.Cells.Item(113, 10).Value = 2749.75
$ws.Cells.Item(113, 11).Value = 1292.7333
$ws.Cells.Item(113, 12).Value = 2749.75
$ws.Cells.Item(113, 13).Value = 877.2666999999999
$ws.Cells.Item(113, 14).Value = -7089.75

$ws.Cells.Item(132, 8).Value = 4289.609
$ws.Cells.Item(132, 9).Value = 3992.611
$ws.Cells.Item(132, 10).Value = 5358.8
$ws.Cells.Item(132, 11).Value = 11977.833
$ws.Cells.Item(132, 12).Value = 16076.4
$ws.Cells.Item(132, 13).Value = -9447.832999999999
$ws.Cells.Item(132, 14).Value = -21136.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 9384.444
$ws.Cells.Item(81, 9).Value = 10432.5
$ws.Cells.Item(81, 10).Value = 1000
$ws.Cells.Item(81, 11).Value = 20865
$ws.Cells.Item(81, 12).Value = 2000
$ws.Cells.Item(81, 13).Value = -19804
$ws.Cells.Item(81, 14).Value = -4122

$ws.Cells.Item(84, 8).Value = 9384.444
$ws.Cells.Item(84, 9).Value = 10432.5
$ws.Cells.Item(84, 10).Value = 1000
$ws.Cells.Item(84, 11).Value = 104325
$ws.Cells.Item(84, 12).Value = 10000
$ws.Cells.Item(84, 13).Value = -99021
$ws.Cells.Item(84, 14).Value = -20608

$ws.Cells.Item(107, 8).Value = 1013.8333
$ws.Cells.Item(107, 9).Value = 969.6667
$ws.Cells.Item(107, 10).Value = 1190.5
$ws.Cells.Item(107, 11).Value = 2909.0001
$ws.Cells.Item(107, 12).Value = 3571.5
$ws.Cells.Item(107, 13).Value = -989.0001000000002
$ws.Cells.Item(107, 14).Value = -7411.5

$ws.Cells.Item(122, 8).Value = 55540
$ws.Cells.Item(122, 9).Value = 60163.043
$ws.Cells.Item(122, 10).Value = 2375
$ws.Cells.Item(122, 11).Value = 180489.129
$ws.Cells.Item(122, 12).Value = 7125
$ws.Cells.Item(122, 13).Value = -178039.129
$ws.Cells.Item(122, 14).Value = -12025

$ws.Cells.Item(126, 8).Value = 158135.38
$ws.Cells.Item(126, 9).Value = 1604.9259
$ws.Cells.Item(126, 11).Value = 4814.7777
$ws.Cells.Item(126, 13).Value = -2344.7777

$ws.Cells.Item(132, 8).Value = 29896.77
$ws.Cells.Item(132, 9).Value = 39617.79
$ws.Cells.Item(132, 10).Value = 3511.1428
$ws.Cells.Item(132, 11).Value = 118853.37
$ws.Cells.Item(132, 12).Value = 10533.4284
$ws.Cells.Item(132, 13).Value = -116323.37
$ws.Cells.Item(132, 14).Value = -15593.4284
